$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Hydrogen / Iron & steel demand (B3)
$ws.Range("B3").Value = 8302.735296243894

# Update Biomass / Non-metallic minerals demand (D6) - tiny floating point correction
$ws.Range("D6").Value = 506.0508721191753
